$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Show your 🤚 hands panther fans 🤚"
$ws.Range("B4").Value = "🧡💙Let's go panthers!💙🧡"
$ws.Range("B8").Value = "@user36 @user37 @user38 @user39 @user40 Today!👇👇👇  "
$ws.Range("B10").Value = "RT @user19 @user42 @user11 no foolin! ⬇️ He’s not playing. "
$ws.Range("B13").Value = "Time for happy hour! Going to Rossi's @user103 https://www.url.com"
$ws.Range("B25").Value = "RT @user19 @user42 @user11 If you here for the rally 🤚🔵 raise your hands for blue 🤚🔵"
$ws.Range("B30").Value = "RT @user49 Get out the Catholic  ✝️ ❤️  vote today🔵Vote Blue!🔵"
$ws.Range("B32").Value = "🚨UPDATE! 🚨 It's all 💩so not happy 😡🥾🐄"
$ws.Range("B76").Value = "🔷Go Blue! ⚾️🔷🔶"
$ws.Range("B88").Value = " 🟠🔵Go orange and blue! 🟠🔵 @user120 Panthers you can do it 🟠🔵"

$ws.Range("B4").Select() | Out-Null
